# Append new thermal-curve measurements (plate2, 20250627, temperature 26)
# for wells A01-A12 and B01-B12 to the bottom of Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$wells  = @("A01","A02","A03","A04","A05","A06","A07","A08","A09","A10","A11","A12",
            "B01","B02","B03","B04","B05","B06","B07","B08","B09","B10","B11","B12")

$lengths = @(5.8730000000000002, 11.772, 14.506, 9.7929999999999993, 11.243, 15.673,
             11.196999999999999, 10.272, 9.048, 9.0579999999999998, 11.081, 13.170999999999999,
             13.327, 4.2489999999999997, 13.641999999999999, 15.699, 14.477, 13.307,
             10.457000000000001, 7.4489999999999998, 9.0749999999999993, 13.901999999999999,
             10.417999999999999, 9.5649999999999995)

$startRow = 98
for ($i = 0; $i -lt $wells.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = 20250627
    $ws.Cells.Item($row, 2).Value = 26
    $ws.Cells.Item($row, 3).Value = "plate2"
    $ws.Cells.Item($row, 4).Value = $wells[$i]
    $ws.Cells.Item($row, 5).Value = $lengths[$i]
}

# Scroll the view down to show the newly added rows and select the next
# empty cell below them, matching where the editor left off.
$ws.Range("C125").Select()

Write-Host "Added $($wells.Length) rows of thermal-curve data (rows $startRow-$($startRow + $wells.Length - 1))"
